# Populate the header row for the vaccinated-faculty sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "FSN"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Age"
$ws.Range("D1").Value = "Phone"
$ws.Range("E1").Value = "Vaccine_Dose"

# Widen the Vaccine_Dose column so the header isn't truncated.
$ws.Columns.Item(5).ColumnWidth = 13.33

# Match the cursor position left behind in the authored workbook.
$ws.Range("F4").Select() | Out-Null
